$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 76: correct the date/time value in column A ---
$ws.Cells.Item(76, 1).Value = 45453.2916666667

# --- Insert new row 77 with data from the R script run ---

# Column A (date): copy formatting (style) from A76 first, so the new
# cell gets the same date/time number format (yyyy-mm-dd hh:mm:ss).
$ws.Cells.Item(76, 1).Copy() | Out-Null
$ws.Cells.Item(77, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(77, 1).Value = 45455.6054282407

$ws.Cells.Item(77, 2).Value = 18000
$ws.Cells.Item(77, 3).Value = 3.27999997138977
$ws.Cells.Item(77, 4).Value = 2.99000000953674
$ws.Cells.Item(77, 5).Value = 2.99000000953674
$ws.Cells.Item(77, 6).Value = 2.99000000953674

# Column G (adj_close) is stored as text "2.99000000953674" (shared string),
# not a number, in the target workbook. Temporarily mark the cell as Text
# so the numeric-looking literal is kept as a string, then restore the
# default (General) formatting/style so no stray style index is left behind.
$gCell = $ws.Cells.Item(77, 7)
$gCell.NumberFormat = "@"
$gCell.Value = "2.99000000953674"
$gCell.Style = "Normal"

# Column H (ticker) is plain text.
$ws.Cells.Item(77, 8).Value = "ESPE.MI"

$excel.CutCopyMode = 0
